# Generate Report for Handoff
# Adds a new handed-off file (c2ca5e9a-52e3-4d71-bb7a-88dbef983b7b) as row 3
# on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newGuid = "c2ca5e9a-52e3-4d71-bb7a-88dbef983b7b"
$newHash = "6ad85fb6f1b26cd40a29226acc57cf9c40d2e5ff"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/f1849a387dd662b1b33066a45ea1adfb898b4882/e2e/$newGuid.md"
$zhAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4058a6bc8dba6961a76ee31424924dc50adf81d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf"
$deAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b394217001699b66a511cd3789008166a5fd8479/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf"

$mdDisplay = "$newGuid.md"
$zhDisplay = "$newGuid.$newHash.zh-cn.xlf"
$deDisplay = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdAddress, "", "", $mdDisplay)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-25 12:46:05"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | ... | Latest Handback DateTime (H) | ...
# Handoff Reason (J)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdAddress, "", "", $mdDisplay)
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhAddress, "", "", $zhDisplay)
$wsZh.Range("E3").Value = "2016-03-25 12:45:59"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("J3").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdAddress, "", "", $mdDisplay)
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deAddress, "", "", $deDisplay)
$wsDe.Range("E3").Value = "2016-03-25 12:46:05"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("J3").Value = "Include"
